# ARS data update: add latest report (2021-02-23 / week 5) and move the
# scratch "helper" calculation block (columns I/J) down below the newly
# appended data, recomputed with fresh numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Duplicate the formatting of the old helper block onto its new
#    location FIRST (before the old cells get cleared below). This is
#    just a format/value carrier - the real content is overwritten in
#    step 3.
# ---------------------------------------------------------------------
$ws.Cells.Item(99, 2).Copy($ws.Cells.Item(169, 2))          # marker cell style

$ws.Range("A149:E154").Copy($ws.Range("A170:E175"))         # weekly data block

$ws.Cells.Item(149, 9).Copy($ws.Cells.Item(169, 9))
$ws.Cells.Item(149, 10).Copy($ws.Cells.Item(169, 10))
$ws.Cells.Item(150, 9).Copy($ws.Cells.Item(170, 9))
$ws.Cells.Item(150, 10).Copy($ws.Cells.Item(170, 10))
$ws.Cells.Item(151, 9).Copy($ws.Cells.Item(171, 9))
$ws.Cells.Item(151, 10).Copy($ws.Cells.Item(171, 10))
$ws.Cells.Item(152, 9).Copy($ws.Cells.Item(172, 9))
$ws.Cells.Item(152, 10).Copy($ws.Cells.Item(172, 10))
$ws.Cells.Item(153, 9).Copy($ws.Cells.Item(173, 9))
$ws.Cells.Item(153, 10).Copy($ws.Cells.Item(173, 10))
$ws.Cells.Item(154, 9).Copy($ws.Cells.Item(174, 9))
$ws.Cells.Item(154, 10).Copy($ws.Cells.Item(174, 10))
$ws.Cells.Item(155, 9).Copy($ws.Cells.Item(175, 9))
$ws.Cells.Item(155, 10).Copy($ws.Cells.Item(175, 10))

$ws.Cells.Item(158, 9).Copy($ws.Cells.Item(178, 9))
$ws.Cells.Item(159, 9).Copy($ws.Cells.Item(179, 9))
$ws.Cells.Item(160, 9).Copy($ws.Cells.Item(180, 9))

# ---------------------------------------------------------------------
# 2. Now remove the old helper block content (columns I/J only - the
#    A:G data columns of rows 149-154/158-160 stay untouched).
# ---------------------------------------------------------------------
$ws.Range("I149:J154").Clear()
$ws.Rows.Item(155).Clear()
$ws.Range("I158").Clear()
$ws.Range("I159").Clear()
$ws.Range("I160").Clear()

# ---------------------------------------------------------------------
# 3. Fill in the real content for the new rows.
# ---------------------------------------------------------------------

# New marker row 169
$ws.Cells.Item(169, 2).Value = "Source: 2021-02-23"

# New weekly data rows 170-175 (2021, week 5)
$ws.Cells.Item(170, 1).Value = 2021
$ws.Cells.Item(170, 2).Value = 5
$ws.Cells.Item(170, 3).Value = "0-4"
$ws.Cells.Item(170, 4).Value = 7819
$ws.Cells.Item(170, 5).Value = 6.7

$ws.Cells.Item(171, 1).Value = 2021
$ws.Cells.Item(171, 2).Value = 5
$ws.Cells.Item(171, 3).Value = "5-14"
$ws.Cells.Item(171, 4).Value = 10664
$ws.Cells.Item(171, 5).Value = 9.8

$ws.Cells.Item(172, 1).Value = 2021
$ws.Cells.Item(172, 2).Value = 5
$ws.Cells.Item(172, 3).Value = "15-34"
$ws.Cells.Item(172, 4).Value = 95972
$ws.Cells.Item(172, 5).Value = 6.8

$ws.Cells.Item(173, 1).Value = 2021
$ws.Cells.Item(173, 2).Value = 5
$ws.Cells.Item(173, 3).Value = "35-59"
$ws.Cells.Item(173, 4).Value = 153555
$ws.Cells.Item(173, 5).Value = 6.5

$ws.Cells.Item(174, 1).Value = 2021
$ws.Cells.Item(174, 2).Value = 5
$ws.Cells.Item(174, 3).Value = "60-79"
$ws.Cells.Item(174, 4).Value = 86730
$ws.Cells.Item(174, 5).Value = 7.3

$ws.Cells.Item(175, 1).Value = 2021
$ws.Cells.Item(175, 2).Value = 5
$ws.Cells.Item(175, 3).Value = ">=80"
$ws.Cells.Item(175, 4).Value = 53318
$ws.Cells.Item(175, 5).Value = 11.3

# New helper block in columns I/J, rows 169-175
$ws.Cells.Item(169, 9).Value = "Länge für Einheit (cm)"
$ws.Cells.Item(169, 10).Formula = "=26.4-6"

$ws.Cells.Item(170, 9).Value = "#:"
$ws.Cells.Item(170, 10).Value = 10

$ws.Cells.Item(171, 9).Value = "# / cm"
$ws.Cells.Item(171, 10).Formula = "=J170 / J169"

$ws.Cells.Item(172, 9).Value = "Achsenabschnitt (cm)"
$ws.Cells.Item(172, 10).Value = 6

$ws.Cells.Item(173, 9).Value = "Achsenabschnitt (#)"
$ws.Cells.Item(173, 10).Value = 5

$ws.Cells.Item(174, 9).Value = "Gemessene Höhe (cm)"
$ws.Cells.Item(174, 10).Value = 9

$ws.Cells.Item(175, 9).Value = "Zahl"
$ws.Cells.Item(175, 10).Formula = "=(J174-J172)*J171+J173"

# Scratch arithmetic moved further down to rows 178-180 (keeping the
# usual two blank-row gap pattern used between weekly blocks).
$ws.Cells.Item(178, 9).Formula = "=50000/55"
$ws.Cells.Item(179, 9).Formula = "=50000 + 909 * 9"
$ws.Cells.Item(180, 9).Formula = "=909*8"

# ---------------------------------------------------------------------
# 4. View-state bookkeeping: scroll the frozen pane and selection down
#    to match the freshly appended rows.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A161"))
$ws.Range("E176").Select()
